$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3419.2683
$ws.Range("J17").Value = 3419.2683
$ws.Range("L17").Value = 10257.8049
$ws.Range("N17").Value = -10593.8049
$ws.Range("H103").Value = 5464146
$ws.Range("I103").Value = 15025202
$ws.Range("J103").Value = 685.7143
$ws.Range("K103").Value = 45075606
$ws.Range("L103").Value = 2057.1429
$ws.Range("M103").Value = -45075020
$ws.Range("N103").Value = -3229.1429
$ws.Range("H132").Value = 3916.3333
$ws.Range("I132").Value = 3721.3
$ws.Range("K132").Value = 11163.9
$ws.Range("M132").Value = -8633.900000000001
$ws.Range("H135").Value = 2066.8262
$ws.Range("I135").Value = 1597
$ws.Range("J135").Value = 7000
$ws.Range("K135").Value = 14373
$ws.Range("L135").Value = 63000
$ws.Range("M135").Value = -11838
$ws.Range("N135").Value = -68070
$ws.Range("H137").Value = 5883199.5
$ws.Range("I137").Value = 872.0909
$ws.Range("J137").Value = 200000000
$ws.Range("K137").Value = 2616.2727
$ws.Range("L137").Value = 600000000
$ws.Range("M137").Value = -66.27269999999999
$ws.Range("N137").Value = -600005100
$ws.Range("H138").Value = 2475.4443
$ws.Range("I138").Value = 1911.6875
$ws.Range("J138").Value = 2926.45
$ws.Range("K138").Value = 5735.0625
$ws.Range("L138").Value = 8779.349999999999
$ws.Range("M138").Value = -595.0625
$ws.Range("N138").Value = -19059.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8368.540999999999
$ws.Range("I32").Value = 9279.416999999999
$ws.Range("J32").Value = 5005.3076
$ws.Range("K32").Value = 9279.416999999999
$ws.Range("L32").Value = 5005.3076
$ws.Range("M32").Value = -8992.416999999999
$ws.Range("N32").Value = -5579.3076
$ws.Range("H61").Value = 29414412
$ws.Range("I61").Value = 29414412
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 29414412
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = -29414200
$ws.Range("M61").ClearContents()
$ws.Range("H74").Value = 10640587
$ws.Range("I74").Value = 15626193
$ws.Range("J74").Value = 4627.933
$ws.Range("K74").Value = 15626193
$ws.Range("L74").Value = 4627.933
$ws.Range("M74").Value = -15625319
$ws.Range("N74").Value = -6375.933
$ws.Range("H77").Value = 10640587
$ws.Range("I77").Value = 15626193
$ws.Range("J77").Value = 4627.933
$ws.Range("K77").Value = 78130965
$ws.Range("L77").Value = 23139.665
$ws.Range("M77").Value = -78126597
$ws.Range("N77").Value = -31875.665
$ws.Range("H97").Value = 6549.5884
$ws.Range("I97").Value = 6883.25
$ws.Range("J97").Value = 1211
$ws.Range("K97").Value = 6883.25
$ws.Range("L97").Value = 1211
$ws.Range("M97").Value = -6387.25
$ws.Range("N97").Value = -2203
$ws.Range("H102").Value = 1202.3334
$ws.Range("I102").Value = 974.4286
$ws.Range("K102").Value = 974.4286
$ws.Range("M102").Value = 647.5714
$ws.Range("H122").Value = 8070.3125
$ws.Range("I122").Value = 10159.417
$ws.Range("J122").Value = 1803
$ws.Range("K122").Value = 30478.251
$ws.Range("L122").Value = 5409
$ws.Range("M122").Value = -28028.251
$ws.Range("N122").Value = -10309
$ws.Range("H136").Value = 29414412
$ws.Range("I136").Value = 29414412
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 88243236
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = -88240686
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 299.66666
$ws.Range("I22").Value = 300.5
$ws.Range("J22").Value = 298
$ws.Range("K22").Value = 300.5
$ws.Range("L22").Value = 298
$ws.Range("M22").Value = -127.5
$ws.Range("N22").Value = -644
$ws.Range("H99").Value = 954.44446
$ws.Range("I99").Value = 936.25
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 936.25
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = 561.75
$ws.Range("N99").Value = -4096

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15879260
$ws.Range("I31").Value = 6795.579
$ws.Range("K31").Value = 6795.579
$ws.Range("M31").Value = -6500.579
$ws.Range("H34").Value = 15879260
$ws.Range("I34").Value = 6795.579
$ws.Range("K34").Value = 6795.579
$ws.Range("M34").Value = -6593.579
$ws.Range("H58").Value = 2436.0952
$ws.Range("I58").Value = 922.1818
$ws.Range("J58").Value = 4101.4
$ws.Range("K58").Value = 922.1818
$ws.Range("L58").Value = 4101.4
$ws.Range("M58").Value = -719.1818
$ws.Range("N58").Value = -4507.4
$ws.Range("H62").Value = 2200
$ws.Range("I62").Value = 2200
$ws.Range("K62").Value = 2200
$ws.Range("M62").Value = -1576
$ws.Range("H65").Value = 2200
$ws.Range("I65").Value = 2200
$ws.Range("K65").Value = 11000
$ws.Range("M65").Value = -7880
$ws.Range("H132").Value = 23816770
$ws.Range("I132").Value = 29419362
$ws.Range("J132").Value = 5753
$ws.Range("K132").Value = 88258086
$ws.Range("L132").Value = 17259
$ws.Range("M132").Value = -88255556
$ws.Range("N132").Value = -22319
$ws.Range("H134").Value = 2517.6428
$ws.Range("I134").Value = 2520.5833
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 7561.749899999999
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -5026.749899999999
$ws.Range("N134").Value = -12570
$ws.Range("H136").Value = 2436.0952
$ws.Range("I136").Value = 922.1818
$ws.Range("J136").Value = 4101.4
$ws.Range("K136").Value = 2766.5454
$ws.Range("L136").Value = 12304.2
$ws.Range("M136").Value = -216.5454
$ws.Range("N136").Value = -17404.2
$ws.Range("H140").Value = 30669.584
$ws.Range("J140").Value = 30669.584
$ws.Range("L140").Value = 30669.584
$ws.Range("N140").Value = -41029.584

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1433.5
$ws.Range("I113").Value = 475.15384
$ws.Range("K113").Value = 1425.46152
$ws.Range("M113").Value = 744.5384799999999
$ws.Range("H132").Value = 804.2143
$ws.Range("I132").Value = 556.75
$ws.Range("K132").Value = 5010.75
$ws.Range("M132").Value = -2480.75
$ws.Range("H134").Value = 3229.8147
$ws.Range("I134").Value = 1779.2106
$ws.Range("J134").Value = 6675
$ws.Range("K134").Value = 5337.6318
$ws.Range("L134").Value = 20025
$ws.Range("M134").Value = -267.6318000000001
$ws.Range("N134").Value = -30165

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 30751.217
$ws.Range("I70").Value = 44451.266
$ws.Range("J70").Value = 5063.625
$ws.Range("K70").Value = 44451.266
$ws.Range("L70").Value = 5063.625
$ws.Range("M70").Value = -44181.266
$ws.Range("N70").Value = -5603.625
$ws.Range("H73").Value = 30751.217
$ws.Range("I73").Value = 44451.266
$ws.Range("J73").Value = 5063.625
$ws.Range("K73").Value = 44451.266
$ws.Range("L73").Value = 5063.625
$ws.Range("M73").Value = -43515.266
$ws.Range("N73").Value = -6935.625
$ws.Range("H122").Value = 2779441
$ws.Range("I122").Value = 3705248.8
$ws.Range("J122").Value = 2017.8334
$ws.Range("K122").Value = 11115746.4
$ws.Range("L122").Value = 6053.5002
$ws.Range("M122").Value = -11113296.4
$ws.Range("N122").Value = -10953.5002
$ws.Range("H126").Value = 3444.9
$ws.Range("I126").Value = 2171.6667
$ws.Range("J126").Value = 4852.1577
$ws.Range("K126").Value = 6515.000100000001
$ws.Range("L126").Value = 14556.4731
$ws.Range("M126").Value = -4045.000100000001
$ws.Range("N126").Value = -19496.4731
$ws.Range("H132").Value = 4727.0557
$ws.Range("I132").Value = 3424.762
$ws.Range("K132").Value = 10274.286
$ws.Range("M132").Value = -7744.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1214
$ws.Range("I22").Value = 406
$ws.Range("J22").Value = 1581.2727
$ws.Range("K22").Value = 406
$ws.Range("L22").Value = 1581.2727
$ws.Range("M22").Value = -111
$ws.Range("N22").Value = -2171.2727
$ws.Range("H27").Value = 1214
$ws.Range("I27").Value = 406
$ws.Range("J27").Value = 1581.2727
$ws.Range("K27").Value = 406
$ws.Range("L27").Value = 1581.2727
$ws.Range("M27").Value = -299
$ws.Range("N27").Value = -1795.2727
$ws.Range("H55").Value = 479.4
$ws.Range("I55").Value = 348.9091
$ws.Range("J55").Value = 638.8889
$ws.Range("K55").Value = 348.9091
$ws.Range("L55").Value = 638.8889
$ws.Range("M55").Value = -175.9091
$ws.Range("N55").Value = -984.8889
$ws.Range("H139").Value = 46439.4
$ws.Range("J139").Value = 47082.668
$ws.Range("L139").Value = 47082.668
$ws.Range("N139").Value = -57362.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 45000
$ws.Range("J118").Value = 45000
$ws.Range("L118").Value = 45000
$ws.Range("N118").Value = -48314
$ws.Range("H126").Value = 10328.777
$ws.Range("I126").Value = 3590.8
$ws.Range("K126").Value = 10772.4
$ws.Range("M126").Value = -8302.400000000001
